$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: a batch of match rows had their home/away data (columns F..V)
# swapped with the adjacent row (column A "Indice" and column E "data_partida"
# stay put - only the match details moved rows).
# ---------------------------------------------------------------------------
function Swap-Rows([int]$r1, [int]$r2) {
    $rng1 = $ws.Range("F$r1" + ":V$r1")
    $rng2 = $ws.Range("F$r2" + ":V$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

Swap-Rows 42 43
Swap-Rows 44 45
Swap-Rows 49 50
Swap-Rows 53 54
Swap-Rows 87 88
Swap-Rows 89 90
Swap-Rows 91 92
Swap-Rows 169 170

# ---------------------------------------------------------------------------
# Part 2: two new match rows were appended at the bottom of the sheet
# (rows 178 and 179), extending the used range from A1:V177 to A1:V179.
# ---------------------------------------------------------------------------

# Clone the formatting of the last existing data row (177) onto the two new
# rows so the Indice/date columns keep their styling (bold border / date
# number format) consistent with the rest of the table.
$src = $ws.Range("A177:V177")
$dst = $ws.Range("A178:V179")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

function Set-Row([int]$r, $vals) {
    $ws.Range("A$r").Value2 = $vals[0]
    $ws.Range("B$r").Value2 = $vals[1]
    $ws.Range("C$r").Value2 = $vals[2]
    $ws.Range("D$r").Value2 = $vals[3]
    $ws.Range("E$r").Value2 = $vals[4]
    $ws.Range("F$r").Value2 = $vals[5]
    $ws.Range("G$r").Value2 = $vals[6]
    $ws.Range("H$r").Value2 = $vals[7]
    $ws.Range("I$r").Value2 = $vals[8]
    $ws.Range("J$r").Value2 = $vals[9]
    $ws.Range("K$r").Value2 = $vals[10]
    $ws.Range("L$r").Value2 = $vals[11]
    $ws.Range("M$r").Value2 = $vals[12]
    $ws.Range("N$r").Value2 = $vals[13]
    $ws.Range("O$r").Value2 = $vals[14]
    $ws.Range("P$r").Value2 = $vals[15]
    $ws.Range("Q$r").Value2 = $vals[16]
    $ws.Range("R$r").Value2 = $vals[17]
    $ws.Range("S$r").Value2 = $vals[18]
    $ws.Range("T$r").Value2 = $vals[19]
    $ws.Range("U$r").Value2 = $vals[20]
    $ws.Range("V$r").Value2 = $vals[21]
}

$row178 = @(
    177, "turkey", "super-lig", "2023-2024", 45298.47916666666,
    "Kayserispor", 1, "Sivasspor", 3,
    1.98, "28/12/2024 19:12",
    2.18, "07/01/2024 11:29",
    3.65, "28/12/2024 19:12",
    3.45, "07/01/2024 11:29",
    3.84, "28/12/2024 19:12",
    3.57, "07/01/2024 11:29",
    "https://www.betexplorer.com/football/turkey/super-lig/kayserispor-sivasspor/Cz0ZXVW5/"
)
Set-Row 178 $row178

$row179 = @(
    178, "turkey", "super-lig", "2023-2024", 45298.58333333334,
    "Istanbulspor AS", 1, "Fenerbahce", 5,
    10.19, "28/12/2024 19:12",
    18.74, "07/01/2024 13:59",
    6.35, "28/12/2024 19:12",
    8.48, "07/01/2024 13:59",
    1.25, "28/12/2024 19:12",
    1.16, "07/01/2024 13:59",
    "https://www.betexplorer.com/football/turkey/super-lig/istanbulspor-as-fenerbahce/jRbwXknC/"
)
Set-Row 179 $row179
